$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 366, shifting the existing rows 366:389
# down to 368:391 (dimension grows from A1:T389 to A1:T391).
$ws.Rows("366:367").Insert()

# --- New row 366: weekly "Primera" quote dated 2022-09-22 (serial 44826) ---
$ws.Range("A366").Value = 4
$ws.Range("B366").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C366").Value = "Los Lagos"
$ws.Range("D366").Value = 44826
$ws.Range("E366").Value = 10
$ws.Range("F366").Value = "Fruta"
$ws.Range("G366").Value = 100102
$ws.Range("H366").Value = "Cítricos"
$ws.Range("I366").Value = 100102006
$ws.Range("J366").Value = "Pomelo"
$ws.Range("K366").Value = "Start Ruby"
$ws.Range("L366").Value = "Primera"
$ws.Range("M366").Value = 120
$ws.Range("N366").Value = 14000
$ws.Range("O366").Value = 15000
$ws.Range("P366").Value = 14500
$ws.Range("Q366").Value = "$/caja 14 kilos empedrada"
$ws.Range("R366").Value = "Región de O'Higgins"
$ws.Range("S366").Value = 1036
$ws.Range("T366").Value = 14

# --- New row 367: weekly "Segunda" quote dated 2022-09-22 (serial 44826) ---
$ws.Range("A367").Value = 4
$ws.Range("B367").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C367").Value = "Los Lagos"
$ws.Range("D367").Value = 44826
$ws.Range("E367").Value = 10
$ws.Range("F367").Value = "Fruta"
$ws.Range("G367").Value = 100102
$ws.Range("H367").Value = "Cítricos"
$ws.Range("I367").Value = 100102006
$ws.Range("J367").Value = "Pomelo"
$ws.Range("K367").Value = "Start Ruby"
$ws.Range("L367").Value = "Segunda"
$ws.Range("M367").Value = 80
$ws.Range("N367").Value = 12000
$ws.Range("O367").Value = 12000
$ws.Range("P367").Value = 12000
$ws.Range("Q367").Value = "$/caja 14 kilos empedrada"
$ws.Range("R367").Value = "Región de O'Higgins"
$ws.Range("S367").Value = 857
$ws.Range("T367").Value = 14
